# Generate Report for Handoff
#
# For the four source files that are still "Ready for handoff"
# (3302dbec, 6a7a4015, 9519626d, e81540ff), regenerating the handoff
# report bumped their Priority from "low" to "ht" and refreshed the
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps
# on both the zh-cn and de-de localization sheets (and on the Overview
# rollup, which shares the de-de timestamp string).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Rows 4-7 are the four "Ready for handoff" files on each language sheet.
$rows = 4..7

foreach ($r in $rows) {
    # Priority column (E): low -> ht
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 5).Value = "ht"

    # Latest Handoff Datetime column (H) on the zh-cn sheet
    $zhcn.Cells.Item($r, 8).Value = "2016-08-22 16:33:29"

    # Latest Handoff Datetime column (H) on the de-de sheet
    $dede.Cells.Item($r, 8).Value = "2016-08-22 16:33:34"

    # Latest HO Xliff Generate Date column (G) on the Overview sheet
    $overview.Cells.Item($r, 7).Value = "2016-08-22 16:33:34"
}
